$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the weight value for sample S1 (B2) - cell becomes empty
$ws.Range("B2").ClearContents()

# Update weight ("wt") values for remaining samples to their corrected
# (fractional) values
$ws.Range("B3").Value = 0.02
$ws.Range("B4").Value = 0.05
$ws.Range("B5").Value = 0.1
$ws.Range("B6").Value = 0.15
$ws.Range("B7").Value = 0.3

# Move the active selection to B8
$ws.Range("B8").Select()
